$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.640.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.670.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +4.69%  "
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000199"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.150.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.473.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.671.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("E28").Value = "  -4.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "533.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
